$wb = $excel.ActiveWorkbook

# The localization status changed from "Ready for handoff" to "In Translation".
# This string is shared across the Overview sheet (columns E & F, row 2) and
# the per-locale detail sheets zh-cn / de-de (column C, row 2) - updating the
# shared text once per occurrence keeps every sheet in sync.
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $v = [string]$cell.Value()
        if ($v -eq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
}

# With the status text now shorter ("In Translation" vs "Ready for handoff"),
# the status columns were narrowed accordingly: E & F on Overview, and C on
# both zh-cn and de-de.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1").ColumnWidth = 12.5
$wsOverview.Range("F1").ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C1").ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C1").ColumnWidth = 12.5
